$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting existing rows (2..29) down to (3..30).
$ws.Rows.Item(2).Insert(-4121)

# The freshly inserted row copies formatting from the header row above it.
# Re-apply the plain data-row formatting (matching the row that got pushed
# down to row 3) to the new row 2 cells.
$ws.Cells.Item(2, 1).Style = $ws.Cells.Item(3, 1).Style
$ws.Cells.Item(2, 1).NumberFormat = $ws.Cells.Item(3, 1).NumberFormat
$ws.Cells.Item(2, 2).ClearFormats()
$ws.Cells.Item(2, 3).ClearFormats()

# Populate the new first data row with the latest fuel price entry
$ws.Cells.Item(2, 1).Value = 45770
$ws.Cells.Item(2, 2).Value = 722.956
$ws.Cells.Item(2, 3).Value = 753.256
